$d = $word.ActiveDocument

$replacements = @(
    @("24×86=", "65×15="),
    @("20×21=", "36×12="),
    @("85×65=", "37×17="),
    @("51×97=", "15×93="),
    @("30×78=", "51×81="),
    @("98×85=", "16×40="),
    @("86×67=", "19×33="),
    @("18×15=", "56×65="),
    @("99×76=", "35×98="),
    @("91×14=", "71×44="),
    @("54×57=", "65×63="),
    @("59×62=", "93×19="),
    @("39×12=", "90×79="),
    @("72×91=", "95×21="),
    @("43×24=", "94×43="),
    @("71×47=", "56×39="),
    @("57×49=", "59×66="),
    @("85×45=", "60×54="),
    @("44×69=", "34×99="),
    @("98×27=", "61×33="),
    @("52×60=", "34×57="),
    @("81×65=", "15×14="),
    @("32×81=", "13×26="),
    @("88×39=", "46×98="),
    @("35×61=", "72×34=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
